# The document's footers (Pearson logo, image2.png) and first-page header
# (BTec logo, image1.jpg) carry inline pictures whose docPr/name metadata
# needs to be swapped: the two Pearson logos become "image1.png" and the
# BTec logo becomes "image2.jpg".
$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer (default) -> footer1.xml -> Pearson logo, was "image2.png"
$footerDefault = $sec.Footers.Item(1)
if ($footerDefault.Exists) {
    for ($i = 1; $i -le $footerDefault.Range.InlineShapes.Count; $i++) {
        $footerDefault.Range.InlineShapes.Item($i).Name = "image1.png"
    }
}

# Footer (first page) -> footer2.xml -> Pearson logo, was "image2.png"
$footerFirst = $sec.Footers.Item(2)
if ($footerFirst.Exists) {
    for ($i = 1; $i -le $footerFirst.Range.InlineShapes.Count; $i++) {
        $footerFirst.Range.InlineShapes.Item($i).Name = "image1.png"
    }
}

# Header (first page) -> header1.xml -> BTec logo, was "image1.jpg"
$headerFirst = $sec.Headers.Item(2)
if ($headerFirst.Exists) {
    for ($i = 1; $i -le $headerFirst.Range.InlineShapes.Count; $i++) {
        $headerFirst.Range.InlineShapes.Item($i).Name = "image2.jpg"
    }
}
